$wb = $excel.ActiveWorkbook
$wsDetalle = $wb.Worksheets.Item("Detalle")
$wsResumen = $wb.Worksheets.Item("Resumen_por_estado")

# --- 1) Five individual status corrections: "En Bodega" (9) -> "Migrado" (4) ---
$singleRows = @(34, 94, 211, 252, 494)
foreach ($r in $singleRows) {
    $wsDetalle.Cells.Item($r, 2).Value = 4
    $wsDetalle.Cells.Item($r, 3).Value = "Migrado"
}

# --- 2) Re-synced data block for rows 2002-2050 (Codigo_Punto, ID_Estado_Migracion_Base, Estado_Migracion_Texto, fecha_ruta) ---
$data = @(
    @(2002, 79332, 2, "Alistamiento", 45934),
    @(2003, 79334, 2, "Alistamiento", 45910),
    @(2004, 79446, 9, "En Bodega", 45881),
    @(2005, 79440, 2, "Alistamiento", 45901),
    @(2006, 80085, 2, "Alistamiento", 45912),
    @(2007, 79166, 9, "En Bodega", 45894),
    @(2008, 79162, 9, "En Bodega", 45890),
    @(2009, 79587, 9, "En Bodega", 45894),
    @(2010, 79669, 2, "Alistamiento", 45959),
    @(2011, 79673, 2, "Alistamiento", 45909),
    @(2012, 79677, 2, "Alistamiento", 45909),
    @(2013, 79962, 9, "En Bodega", 45889),
    @(2014, 80082, 2, "Alistamiento", 45922),
    @(2015, 12543, 2, "Alistamiento", 45931),
    @(2016, 79333, 9, "En Bodega", 45894),
    @(2017, 79435, 9, "En Bodega", 45890),
    @(2018, 79439, 2, "Alistamiento", 45929),
    @(2019, 80084, 2, "Alistamiento", 45936),
    @(2020, 79165, 2, "Alistamiento", 45912),
    @(2021, 79204, 2, "Alistamiento", 45923),
    @(2022, 79586, 2, "Alistamiento", 45898),
    @(2023, 79668, 9, "En Bodega", 45884),
    @(2024, 79590, 2, "Alistamiento", 45898),
    @(2025, 79676, 2, "Alistamiento", 45950),
    @(2026, 79809, 9, "En Bodega", 45880),
    @(2027, 79966, 9, "En Bodega", 45895),
    @(2028, 79433, 2, "Alistamiento", 45960),
    @(2029, 79434, 2, "Alistamiento", 45911),
    @(2030, 79438, 2, "Alistamiento", 45908),
    @(2031, 79442, 2, "Alistamiento", 45898),
    @(2032, 80318, 9, "En Bodega", 45881),
    @(2033, 79167, 2, "Alistamiento", 45904),
    @(2034, 79161, 2, "Alistamiento", 45933),
    @(2035, 79437, 9, "En Bodega", 45895),
    @(2036, 79589, 9, "En Bodega", 45890),
    @(2037, 79672, 2, "Alistamiento", 45959),
    @(2038, 79675, 2, "Alistamiento", 45929),
    @(2039, 79808, 2, "Alistamiento", 45909),
    @(2040, 79965, 9, "En Bodega", 45896),
    @(2041, 79432, 2, "Alistamiento", 45919),
    @(2042, 79336, 2, "Alistamiento", 45910),
    @(2043, 79436, 9, "En Bodega", 45888),
    @(2044, 79441, 2, "Alistamiento", 45936),
    @(2045, 80086, 9, "En Bodega", 45898),
    @(2046, 79164, 2, "Alistamiento", 45901),
    @(2047, 79163, 2, "Alistamiento", 45923),
    @(2048, 79588, 9, "En Bodega", 45884),
    @(2049, 79671, 2, "Alistamiento", 45919),
    @(2050, 79674, 2, "Alistamiento", 45919)
)

foreach ($row in $data) {
    $r = $row[0]
    $wsDetalle.Cells.Item($r, 1).Value = $row[1]
    $wsDetalle.Cells.Item($r, 2).Value = $row[2]
    $wsDetalle.Cells.Item($r, 3).Value = $row[3]
    $wsDetalle.Cells.Item($r, 4).Value = $row[4]
}

# --- 3) Update the summary counts on Resumen_por_estado: "En Bodega" -5, "Migrado" +5 ---
$wsResumen.Cells.Item(3, 2).Value = 300
$wsResumen.Cells.Item(4, 2).Value = 15
